# Add a new "2022-Q1" sheet between "2021-Q4" and "总计", populate it with
# fund holding data, and update the "总计" (totals) summary sheet with a new
# row for 2022-Q1, pushing the existing 2021-Q4 row down.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q4"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Header row (row 1)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row 2 - 159855 银华中证影视主题ETF
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'159855"
$newSheet.Range("C2").Value = "银华中证影视主题ETF"
$newSheet.Range("D2").Value = "'0.96"
$newSheet.Range("E2").Value = "'97.27"
$newSheet.Range("F2").Value = "'4.65"
$newSheet.Range("G2").Value = "'0.0446"
$newSheet.Range("H2").Value = 7

# Data row 3 - 516620 国泰中证影视主题ETF
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'516620"
$newSheet.Range("C3").Value = "国泰中证影视主题ETF"
$newSheet.Range("D3").Value = "'0.33"
$newSheet.Range("E3").Value = "'96.08"
$newSheet.Range("F3").Value = "'4.38"
$newSheet.Range("G3").Value = "'0.0145"
$newSheet.Range("H3").Value = 8

# Formatting to match the style used on the "2021-Q4" sheet: bold, centered,
# top-aligned, thin border around header row and the first (index) column.
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$idxRange = $newSheet.Range("A2:A3")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1

# Keep "2021-Q4" as the active/selected tab, as it was before the edit.
$q4Sheet.Activate()

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row for 2022-Q1 above 2021-Q4.
#    NOTE: must (re)fetch this sheet reference *after* the worksheet was
#    added above, since sheet-index based handles captured earlier would
#    now point at the newly inserted sheet instead.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Re-use the exact formatting already applied to the (now shifted down)
# original index cell so the new row's index cell matches it precisely.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.06

# The pre-existing "2021-Q4" row was shifted down from row 2 to row 3; its
# running index (column A) needs to be bumped from 0 to 1.
$totalSheet.Range("A3").Value = 1
